# Aakash and Madujith time sheet updation
#
# Sheet "06-04-2022" (2nd tab): view/selection tweak only
#   - selection collapses from B2:H12 down to just B2
#   - (scroll position A3 -> B6 is not controllable from this host; left as-is)
#
# Sheet "07-04-2022" (3rd / active tab): rewritten timesheet rows for
# Aakash (row 7) and Madujith (row 11), row-height touch-ups, and a
# cleared cell (D12).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "06-04-2022" — just narrow the selection rectangle to B2
# ---------------------------------------------------------------------
$wsApr06 = $wb.Worksheets.Item("06-04-2022")
$wsApr06.Activate()
$wsApr06.Range("B2").Select()

# ---------------------------------------------------------------------
# "07-04-2022" — Aakash's and Madujith's rows, rewritten
# ---------------------------------------------------------------------
$wsApr07 = $wb.Worksheets.Item("07-04-2022")
$wsApr07.Activate()

# Row 7 — Aakash M
$wsApr07.Range("C7").Value = "Redesigning the home page, Dashboard in Admin page"
$wsApr07.Range("D7").Value = "1)Designed Wireframe for the Public and Approver pages`n2)Discussion about redesigning the wireframe `n3)Session about design pattern (abstract factory)"
$wsApr07.Range("E7").Value = "Wire framing for the UI of HR page"
$wsApr07.Range("F7").Value = "4`n1`n"
$wsApr07.Range("G7").Value = "-`n-`n2`n"
$wsApr07.Rows.Item(7).RowHeight = 101.25

# Row 10 height-only refresh
$wsApr07.Rows.Item(10).RowHeight = 63.75

# Row 11 — Madujith M A
$wsApr07.Range("C11").Value = "Redesinging the Home page wireframe"
$wsApr07.Range("D11").Value = "1).Designed  homepage  wireframe. 2)Explored about Prototype design pattern.3)Disscussion meeting about  redesigning the wireframe 4) Dessign pattern session( abstract factory)"
$wsApr07.Range("F11").Value = "3`n-`n1`n"
$wsApr07.Range("G11").Value = "-`n2`n-`n1"
$wsApr07.Rows.Item(11).RowHeight = 76.5

# Row 12 — clear the stray Discussion note, tidy the row height
$wsApr07.Range("D12").Value = ""
$wsApr07.Rows.Item(12).RowHeight = 25.5

# Final selection / active cell on the active tab
$wsApr07.Range("G7").Select()
